# Update row 8 (year 2025) recurrence metrics in metricas_recorrencia_anual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated raw counts
$ws.Range("C8").Value = 1260
$ws.Range("D8").Value = 204
$ws.Range("E8").Value = 1056

# Recomputed rates based on the updated counts
#  F8 (retention_rate)  = D8 (returning this year) / C7 (total customers prior year) * 100
#  G8 (new_rate)        = E8 (new customers) / C8 (total customers) * 100
#  H8 (returning_rate)  = D8 (returning customers) / C8 (total customers) * 100
$ws.Range("F8").Value = 204 / 2438 * 100
$ws.Range("G8").Value = 1056 / 1260 * 100
$ws.Range("H8").Value = 204 / 1260 * 100
